$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Revisão 2" sheet -----------------------------------
# The workbook opens with "Legenda" as the active sheet (activeTab=1), so a
# bare Worksheets.Add() lands the new sheet immediately before it - i.e.
# between "Revisão 1" and "Legenda", exactly where it belongs.
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "Revisão 2"

# --- 2. Fill it with the second round of functional-requirement answers ----
$data = @(
    @("RFUN3.1","Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim"),
    @("RFUN3.2","Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim"),
    @("RFUN3.3","Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim"),
    @("RFUN4.1","Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim"),
    @("RFUN4.2","Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim"),
    @("RFUN4.3","Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim"),
    @("RFUN4.4","Não","Sim","Sim","Não","Não","Não","Sim","Sim","Sim","Sim")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt $data[$r].Length; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Same centered look as the other answer sheet.
$ws2.Range("A1:K7").HorizontalAlignment = -4108
$ws2.Range("A1:K7").VerticalAlignment = -4108

# Narrower first column, like "Revisão 1" after its own resize below.
$ws2.Columns.Item(1).ColumnWidth = 10.6666666666667

# --- 3. "Revisão 1": tighten the column widths ------------------------------
$ws1 = $wb.Worksheets.Item("Revisão 1")
$ws1.Columns.Item(1).ColumnWidth = 10.6666666666667
$ws1.Range("B1:K1").EntireColumn.ColumnWidth = 8.3333333333334

# --- 4. "Legenda": move the cursor to A5 ------------------------------------
$wsLeg = $wb.Worksheets.Item("Legenda")
[void]$wsLeg.Range("A5").Select()

# --- 5. Leave "Revisão 2" as the active/selected sheet, cursor on L7 --------
$ws2.Activate()
[void]$ws2.Range("L7").Select()
